$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("A28").Value = 111477284
$ws.Range("B28").Value = 77515
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = 'Garnlav'
$ws.Range("G28").Value = 'Alectoria sarmentosa'
$ws.Range("H28").Value = '(Ach.) Ach.'
$ws.Range("J28").Value = ""
$ws.Range("K28").Value = ""
$ws.Range("N28").Value = ""
$ws.Range("P28").Value = 'Bäckslåtten, Jmt'
$ws.Range("Q28").Value = 526737.8549399736
$ws.Range("R28").Value = 7126796.298664714
$ws.Range("S28").Value = 15
$ws.Range("AF28").Value = ""
$ws.Range("AW28").Value = 'Elvira Klang'
$ws.Range("AX28").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 29
$ws.Range("A29").Value = 111475399
$ws.Range("B29").Value = 56398
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = 'Tretåig hackspett'
$ws.Range("G29").Value = 'Picoides tridactylus'
$ws.Range("H29").Value = '(Linnaeus, 1758)'
$ws.Range("K29").Value = ''
$ws.Range("L29").Value = ''
$ws.Range("M29").Value = 'äldre spår'
$ws.Range("N29").Value = ''
$ws.Range("P29").Value = 'Bäckslotten, Jmt'
$ws.Range("Q29").Value = 526448.7861015323
$ws.Range("R29").Value = 7127211.491299792
$ws.Range("S29").Value = 10
$ws.Range("AW29").Value = 'Signe Propst'
$ws.Range("AX29").Value = 'Signe Propst'

# Row 30
$ws.Range("A30").Value = 111477149
$ws.Range("J30").Value = ''
$ws.Range("K30").Value = ''
$ws.Range("N30").Value = ''
$ws.Range("P30").Value = 'bäckslåtten, Jmt'
$ws.Range("Q30").Value = 526773.3997162007
$ws.Range("R30").Value = 7126762.703002418
$ws.Range("S30").Value = 10
$ws.Range("AF30").Value = ''
$ws.Range("AW30").Value = 'Filippa Paperin'
$ws.Range("AX30").Value = 'Filippa Paperin, Karl Soler Kinnerbäck, Tore Dahlberg, Melvin Lewin, Elvira Klang, Elicia Olsson, Jonathan Frendel, Astrid Blomberg, Iris Elmér, Ivar Anderberg, Kai Strömberg, Signe Propst, Elias Blad'

# Row 31
$ws.Range("A31").Value = 111477279
$ws.Range("B31").Value = 78579
$ws.Range("E31").Value = 2081
$ws.Range("F31").Value = 'Skrovellav'
$ws.Range("G31").Value = 'Lobaria scrobiculata'
$ws.Range("H31").Value = '(Scop.) DC.'
$ws.Range("Q31").Value = 526633.8889021962
$ws.Range("R31").Value = 7126891.023735311

# Row 32
$ws.Range("A32").Value = 111477262
$ws.Range("B32").Value = 73696
$ws.Range("E32").Value = 6440
$ws.Range("F32").Value = 'Vitgrynig nållav'
$ws.Range("G32").Value = 'Chaenotheca subroscida'
$ws.Range("H32").Value = '(Eitner) Zahlbr.'
$ws.Range("K32").Value = ""
$ws.Range("L32").Value = ""
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = ""
$ws.Range("P32").Value = 'Bäckslåtten, Jmt'
$ws.Range("Q32").Value = 526874.3871010491
$ws.Range("R32").Value = 7126760.97409881
$ws.Range("S32").Value = 15
$ws.Range("AW32").Value = 'Elvira Klang'
$ws.Range("AX32").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 34
$ws.Range("A34").Value = 111475395
$ws.Range("B34").Value = 56398
$ws.Range("E34").Value = 100109
$ws.Range("F34").Value = 'Tretåig hackspett'
$ws.Range("G34").Value = 'Picoides tridactylus'
$ws.Range("H34").Value = '(Linnaeus, 1758)'
$ws.Range("K34").Value = ''
$ws.Range("L34").Value = ''
$ws.Range("M34").Value = 'äldre spår'
$ws.Range("N34").Value = ''
$ws.Range("P34").Value = 'Bäckslotten, Jmt'
$ws.Range("Q34").Value = 526541.4140344799
$ws.Range("R34").Value = 7127319.634518873
$ws.Range("S34").Value = 10
$ws.Range("AW34").Value = 'Signe Propst'
$ws.Range("AX34").Value = 'Signe Propst'

# Row 35
$ws.Range("A35").Value = 111477277
$ws.Range("K35").Value = ""
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = ""
$ws.Range("P35").Value = 'Bäckslåtten, Jmt'
$ws.Range("Q35").Value = 526643.3394980798
$ws.Range("R35").Value = 7126855.0288987
$ws.Range("S35").Value = 15
$ws.Range("AW35").Value = 'Elvira Klang'
$ws.Range("AX35").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 36
$ws.Range("A36").Value = 111477226
$ws.Range("B36").Value = 73696
$ws.Range("E36").Value = 6440
$ws.Range("F36").Value = 'Vitgrynig nållav'
$ws.Range("G36").Value = 'Chaenotheca subroscida'
$ws.Range("H36").Value = '(Eitner) Zahlbr.'
$ws.Range("Q36").Value = 526776.5518203602
$ws.Range("R36").Value = 7126750.55952362

# Row 37
$ws.Range("A37").Value = 111475415
$ws.Range("B37").Value = 78578
$ws.Range("E37").Value = 6458
$ws.Range("F37").Value = 'Lunglav'
$ws.Range("G37").Value = 'Lobaria pulmonaria'
$ws.Range("H37").Value = '(L.) Hoffm.'
$ws.Range("Q37").Value = 526334.1646390257
$ws.Range("R37").Value = 7127180.955839636

# Row 38
$ws.Range("A38").Value = 111475406
$ws.Range("Q38").Value = 526628.8057909949
$ws.Range("R38").Value = 7126874.89753825

# Row 39
$ws.Range("A39").Value = 111475424
$ws.Range("B39").Value = 77515
$ws.Range("E39").Value = 6425
$ws.Range("F39").Value = 'Garnlav'
$ws.Range("G39").Value = 'Alectoria sarmentosa'
$ws.Range("H39").Value = '(Ach.) Ach.'
$ws.Range("Q39").Value = 526562.7125770835
$ws.Range("R39").Value = 7127372.838116477

# Row 40
$ws.Range("A40").Value = 111475405
$ws.Range("Q40").Value = 526402.8306512056
$ws.Range("R40").Value = 7127140.689653471

# Row 56
$ws.Range("A56").Value = 111481787
$ws.Range("B56").Value = 78578
$ws.Range("D56").Value = 'NT'
$ws.Range("E56").Value = 6458
$ws.Range("F56").Value = 'Lunglav'
$ws.Range("G56").Value = 'Lobaria pulmonaria'
$ws.Range("H56").Value = '(L.) Hoffm.'
$ws.Range("Q56").Value = 526745.1073277664
$ws.Range("R56").Value = 7126863.299787878
$ws.Range("AW56").Value = 'Elias Blad'
$ws.Range("AX56").Value = 'Elias Blad, Astrid Blomberg, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg'

# Row 57
$ws.Range("A57").Value = 111481785
$ws.Range("B57").Value = 78579
$ws.Range("E57").Value = 2081
$ws.Range("F57").Value = 'Skrovellav'
$ws.Range("G57").Value = 'Lobaria scrobiculata'
$ws.Range("H57").Value = '(Scop.) DC.'
$ws.Range("Q57").Value = 526629.4958394679
$ws.Range("R57").Value = 7126895.766960836

# Row 58
$ws.Range("A58").Value = 111481779
$ws.Range("B58").Value = 56398
$ws.Range("E58").Value = 100109
$ws.Range("F58").Value = 'Tretåig hackspett'
$ws.Range("G58").Value = 'Picoides tridactylus'
$ws.Range("H58").Value = '(Linnaeus, 1758)'
$ws.Range("K58").Value = ''
$ws.Range("L58").Value = ''
$ws.Range("M58").Value = 'äldre spår'
$ws.Range("N58").Value = ''
$ws.Range("Q58").Value = 526741.1907235509
$ws.Range("R58").Value = 7126863.265793501

# Row 59
$ws.Range("A59").Value = 111477278
$ws.Range("B59").Value = 96368
$ws.Range("D59").Value = 'LC'
$ws.Range("E59").Value = 221952
$ws.Range("F59").Value = 'Spindelblomster'
$ws.Range("G59").Value = 'Neottia cordata'
$ws.Range("H59").Value = '(L.) Rich.'
$ws.Range("K59").Value = ""
$ws.Range("L59").Value = ""
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = ""
$ws.Range("Q59").Value = 526711.8604537799
$ws.Range("R59").Value = 7126883.005616191
$ws.Range("AW59").Value = 'Elvira Klang'
$ws.Range("AX59").Value = 'Elvira Klang, Tore Dahlberg, Filippa Paperin, Karl Soler Kinnerbäck, Melvin Lewin, Iris Elmér, Signe Propst, Elicia Olsson, Elias Blad, Astrid Blomberg, Jonathan Frendel, Kai Strömberg'

# Row 61
$ws.Range("A61").Value = 111477153
$ws.Range("B61").Value = 78578
$ws.Range("E61").Value = 6458
$ws.Range("F61").Value = 'Lunglav'
$ws.Range("G61").Value = 'Lobaria pulmonaria'
$ws.Range("H61").Value = '(L.) Hoffm.'
$ws.Range("P61").Value = 'bäckslåtten, Jmt'
$ws.Range("Q61").Value = 526793.9223398847
$ws.Range("R61").Value = 7126805.044392107
$ws.Range("AW61").Value = 'Filippa Paperin'
$ws.Range("AX61").Value = 'Filippa Paperin, Karl Soler Kinnerbäck, Tore Dahlberg, Melvin Lewin, Elvira Klang, Elicia Olsson, Jonathan Frendel, Astrid Blomberg, Iris Elmér, Ivar Anderberg, Kai Strömberg, Signe Propst, Elias Blad'

# Row 62
$ws.Range("A62").Value = 111475411
$ws.Range("B62").Value = 78579
$ws.Range("E62").Value = 2081
$ws.Range("F62").Value = 'Skrovellav'
$ws.Range("G62").Value = 'Lobaria scrobiculata'
$ws.Range("H62").Value = '(Scop.) DC.'
$ws.Range("Q62").Value = 526633.8626035146
$ws.Range("R62").Value = 7126894.066088703

# Row 63
$ws.Range("A63").Value = 111475422
$ws.Range("B63").Value = 89423
$ws.Range("E63").Value = 5432
$ws.Range("F63").Value = 'Granticka'
$ws.Range("G63").Value = 'Porodaedalea chrysoloma'
$ws.Range("H63").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q63").Value = 526830.8338072833
$ws.Range("R63").Value = 7126714.5188105

# Row 64
$ws.Range("A64").Value = 111475420
$ws.Range("P64").Value = 'Bäckslotten, Jmt'
$ws.Range("Q64").Value = 526782.4503789117
$ws.Range("R64").Value = 7126672.803025675
$ws.Range("AW64").Value = 'Signe Propst'
$ws.Range("AX64").Value = 'Signe Propst'

# Row 72
$ws.Range("A72").Value = 111475427
$ws.Range("B72").Value = 77515
$ws.Range("E72").Value = 6425
$ws.Range("F72").Value = 'Garnlav'
$ws.Range("G72").Value = 'Alectoria sarmentosa'
$ws.Range("H72").Value = '(Ach.) Ach.'
$ws.Range("Q72").Value = 526432.8127084307
$ws.Range("R72").Value = 7127044.459845332

# Row 73
$ws.Range("A73").Value = 111477156
$ws.Range("B73").Value = 78611
$ws.Range("E73").Value = 6463
$ws.Range("F73").Value = 'Bårdlav'
$ws.Range("G73").Value = 'Nephroma parile'
$ws.Range("H73").Value = '(Ach.) Ach.'
$ws.Range("P73").Value = 'bäckslåtten, Jmt'
$ws.Range("Q73").Value = 526783.9659439438
$ws.Range("R73").Value = 7126798.872483394
$ws.Range("S73").Value = 10
$ws.Range("AW73").Value = 'Filippa Paperin'
$ws.Range("AX73").Value = 'Filippa Paperin, Karl Soler Kinnerbäck, Tore Dahlberg, Melvin Lewin, Elvira Klang, Elicia Olsson, Jonathan Frendel, Astrid Blomberg, Iris Elmér, Ivar Anderberg, Kai Strömberg, Signe Propst, Elias Blad'

# Row 74
$ws.Range("A74").Value = 111477154
$ws.Range("B74").Value = 96348
$ws.Range("D74").Value = 'VU'
$ws.Range("E74").Value = 220787
$ws.Range("F74").Value = 'Knärot'
$ws.Range("G74").Value = 'Goodyera repens'
$ws.Range("H74").Value = '(L.) R. Br.'
$ws.Range("P74").Value = 'bäckslåtten, Jmt'
$ws.Range("Q74").Value = 526784.7834235848
$ws.Range("R74").Value = 7126804.964935203
$ws.Range("S74").Value = 10
$ws.Range("AW74").Value = 'Filippa Paperin'
$ws.Range("AX74").Value = 'Filippa Paperin, Karl Soler Kinnerbäck, Tore Dahlberg, Melvin Lewin, Elvira Klang, Elicia Olsson, Jonathan Frendel, Astrid Blomberg, Iris Elmér, Ivar Anderberg, Kai Strömberg, Signe Propst, Elias Blad'

# Row 75
$ws.Range("A75").Value = 111478257
$ws.Range("B75").Value = 96368
$ws.Range("D75").Value = 'LC'
$ws.Range("E75").Value = 221952
$ws.Range("F75").Value = 'Spindelblomster'
$ws.Range("G75").Value = 'Neottia cordata'
$ws.Range("H75").Value = '(L.) Rich.'
$ws.Range("P75").Value = 'Bäckslåtten, Jmt'
$ws.Range("Q75").Value = 526504.6649119424
$ws.Range("R75").Value = 7127291.938577281
$ws.Range("AW75").Value = 'Elicia Olsson'
$ws.Range("AX75").Value = 'Elicia Olsson, Astrid Blomberg, Elias Blad, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg'

# Row 76
$ws.Range("A76").Value = 111475417
$ws.Range("B76").Value = 78578
$ws.Range("D76").Value = 'NT'
$ws.Range("E76").Value = 6458
$ws.Range("F76").Value = 'Lunglav'
$ws.Range("G76").Value = 'Lobaria pulmonaria'
$ws.Range("H76").Value = '(L.) Hoffm.'
$ws.Range("P76").Value = 'Bäckslotten, Jmt'
$ws.Range("Q76").Value = 526628.5804176665
$ws.Range("R76").Value = 7126900.974891847
$ws.Range("AW76").Value = 'Signe Propst'
$ws.Range("AX76").Value = 'Signe Propst'

# Row 77
$ws.Range("A77").Value = 111481793
$ws.Range("B77").Value = 78612
$ws.Range("D77").Value = 'LC'
$ws.Range("E77").Value = 6464
$ws.Range("F77").Value = 'Luddlav'
$ws.Range("G77").Value = 'Nephroma resupinatum'
$ws.Range("H77").Value = '(L.) Ach.'
$ws.Range("P77").Value = 'Bäckslåtten, Jmt'
$ws.Range("Q77").Value = 526742.0271302022
$ws.Range("R77").Value = 7126867.184979456
$ws.Range("S77").Value = 15
$ws.Range("AW77").Value = 'Elias Blad'
$ws.Range("AX77").Value = 'Elias Blad, Astrid Blomberg, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg'

# Row 78
$ws.Range("A78").Value = 111481780
$ws.Range("B78").Value = 78605
$ws.Range("E78").Value = 6462
$ws.Range("F78").Value = 'Stuplav'
$ws.Range("G78").Value = 'Nephroma bellum'
$ws.Range("H78").Value = '(Spreng.) Tuck.'
$ws.Range("Q78").Value = 526629.4958394679
$ws.Range("R78").Value = 7126895.766960836
$ws.Range("S78").Value = 15
$ws.Range("AW78").Value = 'Elias Blad'
$ws.Range("AX78").Value = 'Elias Blad, Astrid Blomberg, Elicia Olsson, Elvira Klang, Filippa Paperin, Iris Elmér, Ivar Anderberg, Jonathan Frendel, Kai Strömberg, Karl Soler Kinnerbäck, Melvin Lewin, Signe Propst, Tore Dahlberg'

# Row 105
$ws.Range("A105").Value = 111492241
$ws.Range("B105").Value = 89401
$ws.Range("E105").Value = 1108
$ws.Range("F105").Value = 'Harticka'
$ws.Range("G105").Value = 'Pelloporus leporinus'
$ws.Range("H105").Value = '(Fr.) Krieglst.'
$ws.Range("Q105").Value = 526475.7192807253
$ws.Range("R105").Value = 7127216.937948915
$ws.Range("AW105").Value = 'Karl Soler Kinnerbäck'
$ws.Range("AX105").Value = 'Karl Soler Kinnerbäck, Elicia Olsson, Signe Propst, Tore Dahlberg, Melvin Lewin, Elvira Klang, Filippa Paperin, Elias Blad, Iris Elmér, Ivar Anderberg, Kai Strömberg, Jonathan Frendel'

# Row 106
$ws.Range("A106").Value = 111480527
$ws.Range("B106").Value = 76918
$ws.Range("D106").Value = 'NT'
$ws.Range("E106").Value = 6437
$ws.Range("F106").Value = 'Blanksvart spiklav'
$ws.Range("G106").Value = 'Calicium denigratum'
$ws.Range("H106").Value = '(Vain.) Tibell'
$ws.Range("Q106").Value = 526547.1610996595
$ws.Range("R106").Value = 7126854.633392128

# Row 107
$ws.Range("A107").Value = 111480379
$ws.Range("B107").Value = 96348
$ws.Range("D107").Value = 'VU'
$ws.Range("E107").Value = 220787
$ws.Range("F107").Value = 'Knärot'
$ws.Range("G107").Value = 'Goodyera repens'
$ws.Range("H107").Value = '(L.) R. Br.'
$ws.Range("Q107").Value = 526799.075599193
$ws.Range("R107").Value = 7126712.93839974
$ws.Range("AW107").Value = 'Tore Dahlberg'
$ws.Range("AX107").Value = 'Tore Dahlberg, Elvira Klang, Elicia Olsson, Filippa Paperin, Jonathan Frendel, Karl Soler Kinnerbäck, Elias Blad, Signe Propst, Ivar Anderberg, Kai Strömberg, Astrid Blomberg, Melvin Lewin, Iris Elmér'

# Row 132
$ws.Range("A132").Value = 111480173
$ws.Range("B132").Value = 78579
$ws.Range("E132").Value = 2081
$ws.Range("F132").Value = 'Skrovellav'
$ws.Range("G132").Value = 'Lobaria scrobiculata'
$ws.Range("H132").Value = '(Scop.) DC.'
$ws.Range("Q132").Value = 526736.9612626054
$ws.Range("R132").Value = 7126648.50018348

# Row 133
$ws.Range("A133").Value = 111480111
$ws.Range("B133").Value = 56398
$ws.Range("E133").Value = 100109
$ws.Range("F133").Value = 'Tretåig hackspett'
$ws.Range("G133").Value = 'Picoides tridactylus'
$ws.Range("H133").Value = '(Linnaeus, 1758)'
$ws.Range("K133").Value = ''
$ws.Range("L133").Value = ''
$ws.Range("M133").Value = 'äldre spår'
$ws.Range("N133").Value = ''
$ws.Range("Q133").Value = 526775.4342201974
$ws.Range("R133").Value = 7126678.827645465

# Row 134
$ws.Range("A134").Value = 111480065
$ws.Range("B134").Value = 89405
$ws.Range("E134").Value = 1202
$ws.Range("F134").Value = 'Ullticka'
$ws.Range("G134").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H134").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("K134").Value = ""
$ws.Range("L134").Value = ""
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = ""
$ws.Range("Q134").Value = 526618.0915837138
$ws.Range("R134").Value = 7126855.679999291

# Row 142
$ws.Range("A142").Value = 111480130
$ws.Range("B142").Value = 56543
$ws.Range("E142").Value = 103021
$ws.Range("F142").Value = 'Talltita'
$ws.Range("G142").Value = 'Poecile montanus'
$ws.Range("H142").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("Q142").Value = 526894.4666650819
$ws.Range("R142").Value = 7126754.194508768

# Row 143
$ws.Range("A143").Value = 111480175
$ws.Range("B143").Value = 78579
$ws.Range("E143").Value = 2081
$ws.Range("F143").Value = 'Skrovellav'
$ws.Range("G143").Value = 'Lobaria scrobiculata'
$ws.Range("H143").Value = '(Scop.) DC.'
$ws.Range("Q143").Value = 526533.0795019253
$ws.Range("R143").Value = 7127327.385405498
